$d = $word.ActiveDocument

# 1. Bold the first paragraph (heading), including the paragraph mark itself,
#    so both <w:pPr><w:rPr><w:b/></w:rPr></w:pPr> and the run get <w:b/>.
$d.Paragraphs(1).Range.Font.Bold = 1

# 2. Justify (both) the second paragraph (the long description).
$d.Paragraphs(2).Alignment = 3

# 3. Move the "_GoBack" bookmark so it starts at the very beginning of the
#    document (start of paragraph 1) and ends right after the second
#    paragraph (just before the following empty paragraph). Re-adding a
#    bookmark with the same name relocates it, removing the old
#    bookmarkStart/bookmarkEnd pair that used to sit in the spacing
#    paragraph further down.
$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs(3).Range.End
$bookmarkRange = $d.Range($start, $end)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
